$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 8

$ws.Cells.Item($newRow, 1).Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$ws.Cells.Item($newRow, 2).Value = "Beste klant,
Dank voor uw e-mail. Om u beter van dienst te kunnen zijn, heb ik meer specifieke informatie nodig over het product waar u naar informeert. Kunt u ons het productnummer, de naam van het product of enige details geven zodat we het kunnen controleren in ons systeem? 
Met vriendelijke groet,
[Bedrijfsnaam]"
$ws.Cells.Item($newRow, 3).Value = "Kun je nagaan of dit nog leverbaar is?"
$ws.Cells.Item($newRow, 4).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 5).Value = "Productinformatie"
$ws.Cells.Item($newRow, 6).Value = "2025-07-31 21:39:04"
$ws.Cells.Item($newRow, 7).Value = "Ja"
$ws.Cells.Item($newRow, 8).Value = "Nee"
$ws.Cells.Item($newRow, 9).Value = "Ja"
$ws.Cells.Item($newRow, 10).Value = "Nee"

$ws.Rows.Item($newRow).AutoFit()

$wb.Save()
